# Add the new "MemberStatus" sheet after the last existing sheet (PricingStatus)
$wb = $excel.ActiveWorkbook
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "MemberStatus"

# Populate the MemberStatus table
$newSheet.Range("A1").Value = "Id"
$newSheet.Range("B1").Value = "Name"
$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = "Waiting for the task"
$newSheet.Range("A3").Value = 2
$newSheet.Range("B3").Value = "Active"
$newSheet.Range("A4").Value = 3
$newSheet.Range("B4").Value = "On hold"
$newSheet.Range("A5").Value = 4
$newSheet.Range("B5").Value = "Leaver"

# Sheet-local defined names copied in from the source workbook (Calendar/Categories/Customers
# were scoped to the copied sheet there too)
$newSheet.Names.Add("Calendar", "=#REF!") | Out-Null
$newSheet.Names.Add("Categories", "=#REF!") | Out-Null
$newSheet.Names.Add("Customers", "=#REF!") | Out-Null

# Touch every sheet (including the new one) via a cross-workbook reference so Excel records
# the external link back to the source workbook, then clear the helper formulas again.
$sheetNames = @("EmployeePosition","EmploymentStatus","DayType","CustomerStatus","ProjectStatus","PricingStatus","MemberStatus")
$col = 1
foreach ($sn in $sheetNames) {
    $cell = $newSheet.Cells.Item(100, $col)
    $cell.Formula = "='[TimeKeeperStatusesTest.xlsx]$sn'!A1"
    $col = $col + 1
}
$helperRange = $newSheet.Range($newSheet.Cells.Item(100, 1), $newSheet.Cells.Item(100, $sheetNames.Count))
$helperRange.ClearContents() | Out-Null

# Match the author's final selection / active sheet state
$newSheet.Range("B6").Select() | Out-Null
